# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" right after "总计", holding the new
#   quarter's fund-holding data (3 rows).
# - Restyle the oldest quarter sheet ("2020-Q4") header/index cells to match
#   the other quarter sheets (s=2 border/alignment) instead of the one-off
#   "active tab" style it had, and drop it from being the active tab.
# - Prepend a new summary row to "总计" for 2022-Q4 and renumber the index
#   column for the rows that shifted down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Re-style "2020-Q4" (last tab) so its header/index cells match the
#    rest of the quarter sheets (style index 2) instead of the
#    "currently active tab" style (style index 1), and so the tab is no
#    longer marked as the active/selected one.
# ---------------------------------------------------------------------
$wsOld = $wb.Worksheets.Item("2021-Q1")
$wsLast = $wb.Worksheets.Item("2020-Q4")

$wsOld.Range("B1:H1").Copy()
$wsLast.Range("B1:H1").PasteSpecial(-4122)

$wsOld.Range("A2").Copy()
$wsLast.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" sheet right after "总计", by duplicating
#    the "2022-Q1" sheet (same 8-column fund-holding layout/styles) and
#    then overwriting its data.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Copy($null, $wsSummary)
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q4"

# The template ("2022-Q1") only has 2 data rows; 2022-Q4 needs 3, so add
# one more by duplicating row 3's formatting down to row 4 (range-only,
# not whole-row, so the used range doesn't balloon to column XFD).
$wsNew.Range("A3:H3").Copy()
$wsNew.Range("A4:H4").PasteSpecial(-4122)

function Set-TextCell($ws, $addr, $text) {
    # Force the literal string (not an auto-converted number) for
    # numeric-looking values (fund codes, percentages, ...), then drop
    # back to the plain/default style so no stray "@"-format style
    # lingers on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - 华安标普全球石油指数（QDII-LOF）A
Set-TextCell $wsNew "B2" "160416"
Set-TextCell $wsNew "C2" "华安标普全球石油指数（QDII-LOF）A"
Set-TextCell $wsNew "D2" "2.81"
Set-TextCell $wsNew "E2" "93.63"
Set-TextCell $wsNew "F2" "5.18"
Set-TextCell $wsNew "G2" "0.1456"
$wsNew.Range("H2").Value = 4

# Row 3 - 华安标普全球石油指数（QDII-LOF）C
Set-TextCell $wsNew "B3" "014982"
Set-TextCell $wsNew "C3" "华安标普全球石油指数（QDII-LOF）C"
Set-TextCell $wsNew "D3" "0.36"
Set-TextCell $wsNew "E3" "93.63"
Set-TextCell $wsNew "F3" "5.18"
Set-TextCell $wsNew "G3" "0.0186"
$wsNew.Range("H3").Value = 4

# Row 4 - 上投摩根全球新兴市场混合（QDII）
$wsNew.Range("A4").Value = 2
Set-TextCell $wsNew "B4" "378006"
Set-TextCell $wsNew "C4" "上投摩根全球新兴市场混合（QDII）"
Set-TextCell $wsNew "D4" "0.44"
Set-TextCell $wsNew "E4" "86.49"
Set-TextCell $wsNew "F4" "2.66"
Set-TextCell $wsNew "G4" "0.0117"
$wsNew.Range("H4").Value = 6

# ---------------------------------------------------------------------
# 3) Update "总计": insert a new row 2 for 2022-Q4 and bump the index
#    column (A) for all the rows that shifted down by one.
# ---------------------------------------------------------------------
$wsSummary.Rows.Item(2).Insert()
$wsSummary.Range("B2:D2").ClearFormats()

$wsSummary.Range("A3").Copy()
$wsSummary.Range("A2").PasteSpecial(-4122)

$wsSummary.Range("A2").Value = 0
$wsSummary.Range("B2").Value = "2022-Q4"
$wsSummary.Range("C2").Value = 3
$wsSummary.Range("D2").Value = 0.18

for ($r = 3; $r -le 8; $r++) {
    $wsSummary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 4) Leave selection/active tab on the summary sheet (matches the
#    original file, where "总计" is first and nothing else is flagged
#    as the active tab).
# ---------------------------------------------------------------------
$wsSummary.Activate()
$wsSummary.Range("A1").Select()
